$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared by the Overview sheet (B2,C2,B3,C3) and by the
#    per-language sheets' Status column (C2,C3 on zh-cn and de-de).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Populate "Latest Target File" (F) and "Latest Handback File" (G) columns
#    on the zh-cn and de-de sheets, for both data rows (2 and 3), and update
#    the "Latest Handback DateTime" (H) column with a real timestamp.
# ---------------------------------------------------------------------------

$mdUrlRow2 = "https://github.com/OpenLocalizationTest/oltest/blob/9e2f4d693f76102239c2e43c768cdba4b770f1a9/e2e/2f9b3e4b-6c7e-47ae-9d13-f115c80dfcbe.md"
$mdUrlRow3 = "https://github.com/OpenLocalizationTest/oltest/blob/9e2f4d693f76102239c2e43c768cdba4b770f1a9/e2e/857049c3-a339-40ea-86e6-c2e133ef1258.md"

$mdNameRow2 = "2f9b3e4b-6c7e-47ae-9d13-f115c80dfcbe.md"
$mdNameRow3 = "857049c3-a339-40ea-86e6-c2e133ef1258.md"

$xlfUrlZhRow2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d9b777bfc22420acedd49d95fa1cd066342afde/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2f9b3e4b-6c7e-47ae-9d13-f115c80dfcbe.0e437e7a51adfa18b544982b4ef5581f23582fcc.zh-cn.xlf"
$xlfUrlZhRow3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d9b777bfc22420acedd49d95fa1cd066342afde/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/857049c3-a339-40ea-86e6-c2e133ef1258.11a94fdc1e02f6841bdaf3dd78f4e590ba844847.zh-cn.xlf"

$xlfUrlDeRow2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/530973197b39cc3e4148f5f5327995001a6ff49b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2f9b3e4b-6c7e-47ae-9d13-f115c80dfcbe.0e437e7a51adfa18b544982b4ef5581f23582fcc.de-de.xlf"
$xlfUrlDeRow3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/530973197b39cc3e4148f5f5327995001a6ff49b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/857049c3-a339-40ea-86e6-c2e133ef1258.11a94fdc1e02f6841bdaf3dd78f4e590ba844847.de-de.xlf"

$xlfNameZhRow2 = "2f9b3e4b-6c7e-47ae-9d13-f115c80dfcbe.0e437e7a51adfa18b544982b4ef5581f23582fcc.zh-cn.xlf"
$xlfNameZhRow3 = "857049c3-a339-40ea-86e6-c2e133ef1258.11a94fdc1e02f6841bdaf3dd78f4e590ba844847.zh-cn.xlf"

$xlfNameDeRow2 = "2f9b3e4b-6c7e-47ae-9d13-f115c80dfcbe.0e437e7a51adfa18b544982b4ef5581f23582fcc.de-de.xlf"
$xlfNameDeRow3 = "857049c3-a339-40ea-86e6-c2e133ef1258.11a94fdc1e02f6841bdaf3dd78f4e590ba844847.de-de.xlf"

# --- zh-cn sheet ---
$wsZh.Range("F2").Value = $mdNameRow2
$wsZh.Range("F2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $mdUrlRow2, "", "", $mdNameRow2) | Out-Null

$wsZh.Range("G2").Value = $xlfNameZhRow2
$wsZh.Range("G2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $xlfUrlZhRow2, "", "", $xlfNameZhRow2) | Out-Null

$wsZh.Range("F3").Value = $mdNameRow3
$wsZh.Range("F3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $mdUrlRow3, "", "", $mdNameRow3) | Out-Null

$wsZh.Range("G3").Value = $xlfNameZhRow3
$wsZh.Range("G3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $xlfUrlZhRow3, "", "", $xlfNameZhRow3) | Out-Null

$wsZh.Range("H2").Value = "2016-03-13 23:15:13"
$wsZh.Range("H3").Value = "2016-03-13 23:15:13"

# --- de-de sheet ---
$wsDe.Range("F2").Value = $mdNameRow2
$wsDe.Range("F2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $mdUrlRow2, "", "", $mdNameRow2) | Out-Null

$wsDe.Range("G2").Value = $xlfNameDeRow2
$wsDe.Range("G2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $xlfUrlDeRow2, "", "", $xlfNameDeRow2) | Out-Null

$wsDe.Range("F3").Value = $mdNameRow3
$wsDe.Range("F3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $mdUrlRow3, "", "", $mdNameRow3) | Out-Null

$wsDe.Range("G3").Value = $xlfNameDeRow3
$wsDe.Range("G3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $xlfUrlDeRow3, "", "", $xlfNameDeRow3) | Out-Null

$wsDe.Range("H2").Value = "2016-03-13 23:15:19"
$wsDe.Range("H3").Value = "2016-03-13 23:15:19"

$wb.Save()
